# Swap the full data content between row pairs (10,11), (12,13) and (20,22).
# The underlying observations (records) were reordered; every field of one
# row moves to the other row and vice versa, while the row number (and the
# worksheet's per-row position) itself stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All columns that carry data for the affected rows (A:AY used range, but we
# only need to touch the ones that are actually populated in any of the rows
# involved in a swap).
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY")

function Swap-Rows($rowA, $rowB) {
    $valsA = @{}
    $valsB = @{}
    foreach ($col in $cols) {
        $valsA[$col] = $ws.Range($col + $rowA).Value()
        $valsB[$col] = $ws.Range($col + $rowB).Value()
    }
    foreach ($col in $cols) {
        $a = $valsA[$col]
        $b = $valsB[$col]
        # Only touch cells whose value actually changes. This avoids
        # needlessly rewriting (and thereby mangling, e.g. date-like text
        # being re-parsed into a real date serial) cells that hold the same
        # value on both rows, and leaves genuinely-empty/absent cells alone.
        if ($a -ne $b) {
            $ws.Range($col + $rowA).Value = $b
            $ws.Range($col + $rowB).Value = $a
        }
    }
}

Swap-Rows 10 11
Swap-Rows 12 13
Swap-Rows 20 22
